$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2026-02-06 Friday" "2026-02-07 Saturday"

Replace-Text "621×9=5589" "750×6=4500"
Replace-Text "955×7=6685" "844×5=4220"
Replace-Text "297×2=594" "464×3=1392"
Replace-Text "692×3=2076" "374×6=2244"
Replace-Text "253×2=506" "668×2=1336"

Replace-Text "263×8=2104" "897×2=1794"
Replace-Text "809×9=7281" "306×3=918"
Replace-Text "277×6=1662" "170×5=850"
Replace-Text "291×8=2328" "754×6=4524"
Replace-Text "295×5=1475" "831×3=2493"

Replace-Text "323×6=1938" "565×4=2260"
Replace-Text "784×5=3920" "273×4=1092"
Replace-Text "736×7=5152" "492×2=984"
Replace-Text "623×9=5607" "582×4=2328"
Replace-Text "687×3=2061" "742×8=5936"

Replace-Text "342×2=684" "819×7=5733"
Replace-Text "781×8=6248" "354×3=1062"
Replace-Text "979×7=6853" "567×5=2835"
Replace-Text "106×5=530" "920×3=2760"
Replace-Text "387×6=2322" "117×5=585"

Replace-Text "352×5=1760" "173×2=346"
Replace-Text "526×6=3156" "224×9=2016"
Replace-Text "502×4=2008" "128×9=1152"
Replace-Text "622×4=2488" "173×5=865"
Replace-Text "703×9=6327" "474×2=948"
